# Aspekt 1, Mixing E3/E4
# Add a new logged activity row (27) to the "Geleistete Arbeiten" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New entry: 04.06.2013 (serial 41429), "Eclipse 4.3 Kepler + Migrationsmöglichkeiten analysiert", 3 hours
$ws.Cells.Item(27, 1).Value2 = 41429
$ws.Cells.Item(27, 2).Value2 = "Eclipse 4.3 Kepler + Migrationsmöglichkeiten analysiert"
$ws.Cells.Item(27, 3).Value2 = 3

# Match the date formatting/style used by the other date cells in column A
# (copy format from the preceding row instead of assigning a NumberFormat
# string, so it reuses the existing style/numFmt rather than creating a new one).
$ws.Cells.Item(26, 1).Copy()
$ws.Cells.Item(27, 1).PasteSpecial(-4122)

# The SUM(C2:C33) total in row 34 recalculates automatically (99 -> 102).

# Update the active selection to reflect where the edit was made.
$ws.Range("B27").Select()
